$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("temperature vs brightness")

# Duplicate rows 6-14 (which pull the per-brightness "temperature rise" figures
# via cross-sheet formulas) into rows 22-30, but referencing the already
# computed A/C columns directly on this sheet, with new "&"/"\\" label
# columns (B/D) flanking the duplicated value, per report automation using
# octave-generated graphs.
for ($i = 0; $i -le 8; $i++) {
    $srcRow = 6 + $i
    $dstRow = 22 + $i

    $ws.Range("A$dstRow").Formula = "=A$srcRow"
    $ws.Range("B$dstRow").Value = "&"
    $ws.Range("C$dstRow").Formula = "=C$srcRow"
    $ws.Range("D$dstRow").Value = "\\"
}
